# Adds the new "tbl_estoque" table (Tabela79) to the "oficina" worksheet,
# mirroring the structure of the existing "tbl_peca" table (Tabela7) but
# with stock/estoque related rows (oleo de motor, filtro de ar, filtro
# combustivel) placed in columns L:O (rows 15-19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "tbl_peca" block (G15:J19) onto the
# new block (L15:O19) so the new table visually matches its sibling.
$srcFormat = $ws.Range("G15:J19")
$dstFormat = $ws.Range("L15:O19")
$srcFormat.Copy()
$dstFormat.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Title band above the table
$ws.Range("L15").Value = "tbl_estoque"

# Header row
$ws.Range("L16").Value = "cod_peca"
$ws.Range("M16").Value = "nome"
$ws.Range("N16").Value = "descricao"
$ws.Range("O16").Value = "cod_servico"

# Data rows
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = "óleo de motor"
$ws.Range("O17").Value = 1

$ws.Range("L18").Value = 2
$ws.Range("M18").Value = "filtro de ar"

$ws.Range("L19").Value = 3
$ws.Range("M19").Value = "filtro combustível"

# Turn the range into a proper Excel Table (ListObject), matching the
# other small lookup tables already on the sheet.
$lo = $ws.ListObjects.Add(1, $ws.Range("L16:O19"), [Type]::Missing, 1)
$lo.Name = "Tabela79"

# Approximate the column width auto-fit that Excel performs once the new
# (wider) text is typed into columns L:O.
$ws.Columns("L:L").ColumnWidth = 15
$ws.Columns("M:M").ColumnWidth = 17
$ws.Columns("N:N").ColumnWidth = 14.833333333333334
$ws.Columns("O:O").ColumnWidth = 17.333333333333332

# Leave the same selection state (whole L:O columns) that is present in
# the saved workbook.
$ws.Columns("L:O").Select() | Out-Null

Write-Output "tbl_estoque table added"
